$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the pass/fail condition text in rows 6 and 7 ---
# Row 6 used to be the "pass" row (>= 80 / сдал итоговый экзамен); it becomes the "fail" row.
# Row 7 used to be the "fail" row (< 80 / завалил итоговый экзамен); it becomes the "pass" row
# with an updated, more complete condition.
$ws.Range("B6").Value = "< 80"
$ws.Range("C6").Value = "завалил итоговый экзамен"
$ws.Range("B7").Value = ">= 80 and <= 100"
$ws.Range("C7").Value = "сдал итоговый экзамен"

# --- Header row is a touch taller ---
$ws.Rows.Item(1).RowHeight = 19.5

# --- A2:B5 and A6:A7 adopt the same (theme-colour) font style as the header instead
#     of the explicit black-RGB font, so copy the header cell's formatting onto them.
#     (B6/B7 keep their own "left aligned with border" style.) ---
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A2:B5").PasteSpecial(-4122) | Out-Null
$ws.Range("A6:A7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

Write-Output "done"
